# Update cryptos list: price (D) and 1h volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text so numeric-looking strings (e.g. "584.17")
# are stored as literal text instead of being coerced to a Number - matches
# the original inline-string cells. ClearFormats() afterwards drops the
# temporary "@" number-format style again so no stray cell style is left
# behind (cells keep style 0, same as before the edit).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.114.47'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.472.28'
$ws.Range('E3').Value = '  -2.87%  '
$ws.Range('D5').Value = '584.17'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').Value = '168.51'
$ws.Range('E6').Value = '  -3.01%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -2.29%  '
$ws.Range('D9').Value = '2.471.90'
$ws.Range('E9').Value = '  -2.82%  '
$ws.Range('E10').Value = '  -2.97%  '
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('E12').Value = '  -2.39%  '
$ws.Range('D13').Value = '0.330'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '25.63'
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('D15').Value = '2.917.05'
$ws.Range('D16').Value = '66.812.13'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('E17').Value = '  -4.96%  '
$ws.Range('D18').Value = '2.454.60'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').Value = '11.10'
$ws.Range('E19').Value = '  -6.04%  '
$ws.Range('E20').Value = '  -5.48%  '
$ws.Range('D21').Value = '354.19'
$ws.Range('E21').Value = '  -4.23%  '
$ws.Range('D22').Value = '4.03'
$ws.Range('E22').Value = '  -3.18%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '68.98'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('E25').Value = '  -7.33%  '
$ws.Range('D26').Value = '1.84'
$ws.Range('E26').Value = '  -5.17%  '
$ws.Range('D27').Value = '9.26'
$ws.Range('E27').Value = '  -7.26%  '
$ws.Range('E28').Value = '  -57.86%  '
$ws.Range('D29').Value = '2.591.44'
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('E30').Value = '  -7.08%  '
$ws.Range('D31').Value = '515.68'
$ws.Range('E31').Value = '  -5.20%  '
$ws.Range('D32').Value = '7.74'
$ws.Range('E32').Value = '  -8.78%  '
$ws.Range('E33').Value = '  -6.32%  '
$ws.Range('E34').Value = '  -5.50%  '
$ws.Range('E36').Value = '  -9.56%  '
$ws.Range('D37').Value = '158.85'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('E39').Value = '  -4.25%  '
$ws.Range('E40').Value = '  -6.54%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '4.84'
$ws.Range('E42').Value = '  -6.49%  '
$ws.Range('E43').Value = '  -6.33%  '
$ws.Range('E44').Value = '  -7.02%  '
$ws.Range('E45').Value = '  -7.65%  '
$ws.Range('D46').Value = '38.84'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Value = '140.29'
$ws.Range('E47').Value = '  -5.16%  '
$ws.Range('E48').Value = '  -7.06%  '
$ws.Range('D49').Value = '0.516'
$ws.Range('E49').Value = '  -6.93%  '
$ws.Range('E50').Value = '  -11.80%  '
$ws.Range('E51').Value = '  -7.34%  '

$priceRange.ClearFormats()
